$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension implicitly by writing data; set new shared string + full data grid

# Row 2
$ws.Range("A2").Value = "M1"
$ws.Range("B2").Value = "Ebi3"
$ws.Range("C2").Value = "Il27ra"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.392984666666666
$ws.Range("H2").Value = 16.178954
$ws.Range("I2").Value = 0.3709566379599203
$ws.Range("J2").Value = 0.3709566379599202
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.899454
$ws.Range("N2").Value = 3.798908
$ws.Range("O2").Value = 0.4692036141039673
$ws.Range("P2").Value = 0.3739095189885477
$ws.Range("Q2").Value = 10.24372629703867
$ws.Range("R2").Value = 61.46235778223199
$ws.Range("S2").Value = 0.1740541952066516
$ws.Range("T2").Value = 0.1387042180652026

# Row 3
$ws.Range("A3").Value = "M1"
$ws.Range("B3").Value = "Ebi3"
$ws.Range("C3").Value = "Il27ra"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.392984666666666
$ws.Range("H3").Value = 16.178954
$ws.Range("I3").Value = 0.3709566379599203
$ws.Range("J3").Value = 0.3709566379599202
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.5130936666666667
$ws.Range("N3").Value = 1.539281
$ws.Range("O3").Value = 0.1267445290982863
$ws.Range("P3").Value = 0.1515045424364609
$ws.Range("Q3").Value = 2.767106276897111
$ws.Range("R3").Value = 24.90395649207399
$ws.Range("S3").Value = 0.04701672439411357
$ws.Range("T3").Value = 0.05620161569788559

# Row 4
$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Ebi3"
$ws.Range("C4").Value = "Il27ra"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.392984666666666
$ws.Range("H4").Value = 16.178954
$ws.Range("I4").Value = 0.3709566379599203
$ws.Range("J4").Value = 0.3709566379599202
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1162836666666667
$ws.Range("N4").Value = 0.348851
$ws.Range("O4").Value = 0.02872442115537467
$ws.Range("P4").Value = 0.03433584324986914
$ws.Range("Q4").Value = 0.6271160313171111
$ws.Range("R4").Value = 5.644044281854
$ws.Range("S4").Value = 0.0106555146991426
$ws.Range("T4").Value = 0.01273710897349028

# Row 5
$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = "Ebi3"
$ws.Range("C5").Value = "Il27ra"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.392984666666666
$ws.Range("H5").Value = 16.178954
$ws.Range("I5").Value = 0.3709566379599203
$ws.Range("J5").Value = 0.3709566379599202
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3438503333333333
$ws.Range("N5").Value = 1.031551
$ws.Range("O5").Value = 0.08493799750394264
$ws.Range("P5").Value = 0.1015309500051476
$ws.Range("Q5").Value = 1.854379575294889
$ws.Range("R5").Value = 16.689416177654
$ws.Range("S5").Value = 0.03150831398911066
$ws.Range("T5").Value = 0.03766357986278632

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Ebi3"
$ws.Range("C6").Value = "Il27ra"
$ws.Range("D6").Value = "Neutro"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.392984666666666
$ws.Range("H6").Value = 16.178954
$ws.Range("I6").Value = 0.3709566379599203
$ws.Range("J6").Value = 0.3709566379599202
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.090236333333333
$ws.Range("N6").Value = 3.270709
$ws.Range("O6").Value = 0.2693104585988698
$ws.Range("P6").Value = 0.3219212544608909
$ws.Range("Q6").Value = 5.879627828709555
$ws.Range("R6").Value = 52.91665045838599
$ws.Range("S6").Value = 0.09990250228928105
$ws.Range("T6").Value = 0.1194188262426521

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Ebi3"
$ws.Range("C7").Value = "Il27ra"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.392984666666666
$ws.Range("H7").Value = 16.178954
$ws.Range("I7").Value = 0.3709566379599203
$ws.Range("J7").Value = 0.3709566379599202
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.5
$ws.Range("M7").Value = 0.08533300000000001
$ws.Range("N7").Value = 0.170666
$ws.Range("O7").Value = 0.02107897953955918
$ws.Range("P7").Value = 0.01679789085908358
$ws.Range("Q7").Value = 0.4601995605606666
$ws.Range("R7").Value = 2.761197363364
$ws.Range("S7").Value = 0.007819387381620822
$ws.Range("T7").Value = 0.006231289117903323

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Ebi3"
$ws.Range("C8").Value = "Il27ra"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 9.145061333333333
$ws.Range("H8").Value = 27.435184
$ws.Range("I8").Value = 0.6290433620400798
$ws.Range("J8").Value = 0.6290433620400798
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.899454
$ws.Range("N8").Value = 3.798908
$ws.Range("O8").Value = 0.4692036141039673
$ws.Range("P8").Value = 0.3739095189885477
$ws.Range("Q8").Value = 17.37062332984533
$ws.Range("R8").Value = 104.223739979072
$ws.Range("S8").Value = 0.2951494188973158
$ws.Range("T8").Value = 0.2352053009233451

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Ebi3"
$ws.Range("C9").Value = "Il27ra"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 9.145061333333333
$ws.Range("H9").Value = 27.435184
$ws.Range("I9").Value = 0.6290433620400798
$ws.Range("J9").Value = 0.6290433620400798
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.5130936666666667
$ws.Range("N9").Value = 1.539281
$ws.Range("O9").Value = 0.1267445290982863
$ws.Range("P9").Value = 0.1515045424364609
$ws.Range("Q9").Value = 4.692273051411555
$ws.Range("R9").Value = 42.230457462704
$ws.Range("S9").Value = 0.07972780470417273
$ws.Range("T9").Value = 0.09530292673857531

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Ebi3"
$ws.Range("C10").Value = "Il27ra"
$ws.Range("D10").Value = "M1"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 9.145061333333333
$ws.Range("H10").Value = 27.435184
$ws.Range("I10").Value = 0.6290433620400798
$ws.Range("J10").Value = 0.6290433620400798
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.1162836666666667
$ws.Range("N10").Value = 0.348851
$ws.Range("O10").Value = 0.02872442115537467
$ws.Range("P10").Value = 0.03433584324986914
$ws.Range("Q10").Value = 1.063421263731555
$ws.Range("R10").Value = 9.570791373584001
$ws.Range("S10").Value = 0.01806890645623208
$ws.Range("T10").Value = 0.02159873427637887

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Ebi3"
$ws.Range("C11").Value = "Il27ra"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 9.145061333333333
$ws.Range("H11").Value = 27.435184
$ws.Range("I11").Value = 0.6290433620400798
$ws.Range("J11").Value = 0.6290433620400798
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.3438503333333333
$ws.Range("N11").Value = 1.031551
$ws.Range("O11").Value = 0.08493799750394264
$ws.Range("P11").Value = 0.1015309500051476
$ws.Range("Q11").Value = 3.144532387820444
$ws.Range("R11").Value = 28.300791490384
$ws.Range("S11").Value = 0.05342968351483199
$ws.Range("T11").Value = 0.06386737014236134

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Ebi3"
$ws.Range("C12").Value = "Il27ra"
$ws.Range("D12").Value = "Neutro"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 9.145061333333333
$ws.Range("H12").Value = 27.435184
$ws.Range("I12").Value = 0.6290433620400798
$ws.Range("J12").Value = 0.6290433620400798
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.090236333333333
$ws.Range("N12").Value = 3.270709
$ws.Range("O12").Value = 0.2693104585988698
$ws.Range("P12").Value = 0.3219212544608909
$ws.Range("Q12").Value = 9.970278136161777
$ws.Range("R12").Value = 89.732503225456
$ws.Range("S12").Value = 0.1694079563095888
$ws.Range("T12").Value = 0.2025024282182388

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Ebi3"
$ws.Range("C13").Value = "Il27ra"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 9.145061333333333
$ws.Range("H13").Value = 27.435184
$ws.Range("I13").Value = 0.6290433620400798
$ws.Range("J13").Value = 0.6290433620400798
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.5
$ws.Range("M13").Value = 0.08533300000000001
$ws.Range("N13").Value = 0.170666
$ws.Range("O13").Value = 0.02107897953955918
$ws.Range("P13").Value = 0.01679789085908358
$ws.Range("Q13").Value = 0.7803755187573334
$ws.Range("R13").Value = 4.682253112544
$ws.Range("S13").Value = 0.01325959215793836
$ws.Range("T13").Value = 0.01056660174118026
